$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.597.22"
$ws.Range("E2").Value = "  -2.48%  "
$ws.Range("D3").Value = "1.664.71"
$ws.Range("E3").Value = "  -3.95%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.509"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.65%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0877"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("D12").Value = "1.901.24"
$ws.Range("E12").Value = "  -3.89%  "
$ws.Range("D13").Value = "1.700.14"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.21%  "
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("D17").Value = "27.600.74"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "241.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.94%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("E23").Value = "  -3.28%  "
$ws.Range("E24").Value = "  -3.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").Value = "1.457.04"
$ws.Range("E33").Value = "  -3.22%  "
$ws.Range("E34").Value = "  -4.37%  "
$ws.Range("E35").Value = "  -3.89%  "
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("E38").Value = "  -4.93%  "
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "69.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("E41").Value = "  -4.58%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.18%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.52%  "
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").Value = "1.808.92"
$ws.Range("E46").Value = "  -3.86%  "
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("E49").Value = "  -5.24%  "
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.40%  "
